$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tree 1")
$ws.Rows.Item(2).Insert()
